# "remove frames from arch.png"
#
# The arch.png diagram on both slides is built from native PowerPoint
# shapes stacked on top of the picture. A handful of text-label boxes
# had a thin frame (outline) drawn around them; this edit removes those
# frames, nudges one label back into alignment, tidies the connector
# "glue" metadata, and drops slide 2's explicit white background so it
# goes back to inheriting the layout/master background.

$p = $ppt.ActivePresentation

# EMU-exact offset fix for the "Match Orders" textbox (id 30) on both
# slides. Shape.Left/.Top on a shape nested in (flattened) groups is
# converted straight to EMU on write, so feed it the exact target EMU
# value expressed in points to land on the exact OOXML offset.
$matchOrdersTargetOffX = 2266618
$matchOrdersTargetLeftPt = $matchOrdersTargetOffX / 12700.0

for ($slideIdx = 1; $slideIdx -le $p.Slides.Count; $slideIdx++) {
    $s = $p.Slides.Item($slideIdx)
    $top = $s.Shapes.Item(1)
    $items = $top.GroupItems

    # --- Connectors: restore the "glued" connector lock marker ------
    # Only slide 1's three loose (previously-unlocked) connectors gain
    # the lock; the rest of the deck is left untouched.
    if ($slideIdx -eq 1) {
        $connectorNames = @("直接箭头连接符 20", "连接符: 肘形 35", "连接符: 肘形 38")
        foreach ($cname in $connectorNames) {
            $conn = $items.Item($cname)
            $conn.LockAspectRatio = $true
        }
    }

    # --- Reposition the "Match Orders" text box ----------------------
    $matchOrders = $items.Item("文本框 29")
    $matchOrders.Left = $matchOrdersTargetLeftPt

    # --- Slide 2 only: drop per-shape frame outlines & background ---
    if ($slideIdx -eq 2) {
        $s.FollowMasterBackground = $true

        $labelNames = @(
            "文本框 11",
            "文本框 18",
            "文本框 22",
            "文本框 28",
            "文本框 29",
            "文本框 37",
            "文本框 47",
            "文本框 48"
        )
        foreach ($lname in $labelNames) {
            $shp = $items.Item($lname)
            $shp.Line.Visible = $false
        }
    }
}
